$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 68.666664
$ws.Range("I8").Value = 68.666664
$ws.Range("K8").Value = 205.999992
$ws.Range("M8").Value = -66.99999199999999

$ws.Range("H17").Value = 5639271.5
$ws.Range("J17").Value = 5639271.5
$ws.Range("L17").Value = 16917814.5
$ws.Range("N17").Value = -16918150.5

$ws.Range("H40").Value = 3128.1428
$ws.Range("J40").Value = 3284.875
$ws.Range("L40").Value = 3284.875
$ws.Range("N40").Value = -3634.875

$ws.Range("H55").Value = 909
$ws.Range("I55").Value = 399.8
$ws.Range("J55").Value = 1333.3334
$ws.Range("K55").Value = 399.8
$ws.Range("L55").Value = 1333.3334
$ws.Range("M55").Value = -185.8
$ws.Range("N55").Value = -1761.3334

$ws.Range("H111").Value = 2509.7058
$ws.Range("I111").Value = 2332.1538
$ws.Range("K111").Value = 6996.4614
$ws.Range("M111").Value = -3929.4614

$ws.Range("H131").Value = 7398.1
$ws.Range("I131").Value = 7398.1
$ws.Range("K131").Value = 22194.3
$ws.Range("M131").Value = -17154.3

$ws.Range("H133").Value = 124993.75
$ws.Range("J133").Value = 124993.75
$ws.Range("L133").Value = 124993.75
$ws.Range("N133").Value = -135113.75

$ws.Range("H135").Value = 5283.727
$ws.Range("I135").Value = 5569.5557
$ws.Range("K135").Value = 50126.0013
$ws.Range("M135").Value = -47591.0013

$ws.Range("H138").Value = 6108.104
$ws.Range("I138").Value = 13234.333
$ws.Range("J138").Value = 2868.9092
$ws.Range("K138").Value = 39702.999
$ws.Range("L138").Value = 8606.7276
$ws.Range("M138").Value = -34562.999
$ws.Range("N138").Value = -18886.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 200670.02
$ws.Range("I32").Value = 214920.36
$ws.Range("J32").Value = 33228.5
$ws.Range("K32").Value = 214920.36
$ws.Range("L32").Value = 33228.5
$ws.Range("M32").Value = -214633.36
$ws.Range("N32").Value = -33802.5

$ws.Range("H45").Value = 114186.78
$ws.Range("I45").Value = 128272.75
$ws.Range("K45").Value = 128272.75
$ws.Range("M45").Value = -127895.75

$ws.Range("H74").Value = 329214.94
$ws.Range("I74").Value = 1340.3948
$ws.Range("J74").Value = 744522.7
$ws.Range("K74").Value = 1340.3948
$ws.Range("L74").Value = 744522.7
$ws.Range("M74").Value = -466.3948
$ws.Range("N74").Value = -746270.7

$ws.Range("H77").Value = 329214.94
$ws.Range("I77").Value = 1340.3948
$ws.Range("J77").Value = 744522.7
$ws.Range("K77").Value = 6701.974
$ws.Range("L77").Value = 3722613.5
$ws.Range("M77").Value = -2333.974
$ws.Range("N77").Value = -3731349.5

$ws.Range("H102").Value = 2413.0476
$ws.Range("I102").Value = 2367.0527
$ws.Range("K102").Value = 2367.0527
$ws.Range("M102").Value = -745.0527000000002

$ws.Range("H132").Value = 1454.069
$ws.Range("I132").Value = 794.34784
$ws.Range("K132").Value = 2383.04352
$ws.Range("M132").Value = 146.9564799999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8370.333000000001
$ws.Range("I105").Value = 9516.286
$ws.Range("K105").Value = 9516.286
$ws.Range("M105").Value = -7769.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2650.6
$ws.Range("I16").Value = 1584.45
$ws.Range("K16").Value = 1584.45
$ws.Range("M16").Value = -1297.45

$ws.Range("H31").Value = 2937.5417
$ws.Range("I31").Value = 2223.4
$ws.Range("J31").Value = 3020.5813
$ws.Range("K31").Value = 2223.4
$ws.Range("L31").Value = 3020.5813
$ws.Range("M31").Value = -1928.4
$ws.Range("N31").Value = -3610.5813

$ws.Range("H34").Value = 2937.5417
$ws.Range("I34").Value = 2223.4
$ws.Range("J34").Value = 3020.5813
$ws.Range("K34").Value = 2223.4
$ws.Range("L34").Value = 3020.5813
$ws.Range("M34").Value = -2021.4
$ws.Range("N34").Value = -3424.5813

$ws.Range("H96").Value = 31000
$ws.Range("J96").Value = 31000
$ws.Range("L96").Value = 31000
$ws.Range("N96").Value = -36492

$ws.Range("H113").Value = 2650.6
$ws.Range("I113").Value = 1584.45
$ws.Range("K113").Value = 1584.45
$ws.Range("M113").Value = 585.55

$ws.Range("H122").Value = 3017.9546
$ws.Range("I122").Value = 2516.6667
$ws.Range("K122").Value = 7550.000100000001
$ws.Range("M122").Value = -5100.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 373.47058
$ws.Range("I2").Value = 218.72
$ws.Range("K2").Value = 1312.32
$ws.Range("M2").Value = -1199.32

$ws.Range("H108").Value = 333333630
$ws.Range("I108").Value = 333333630
$ws.Range("K108").Value = 1000000890
$ws.Range("M108").Value = -999998010

$ws.Range("H113").Value = 1060.7273
$ws.Range("J113").Value = 972.25
$ws.Range("L113").Value = 2916.75
$ws.Range("N113").Value = -7256.75

$ws.Range("H121").Value = 1799.2
$ws.Range("J121").Value = 1799.2
$ws.Range("L121").Value = 5397.6
$ws.Range("N121").Value = -8017.6

$ws.Range("H131").Value = 10154148
$ws.Range("J131").Value = 8627614
$ws.Range("L131").Value = 25882842
$ws.Range("N131").Value = -25892922

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 21206054
$ws.Range("I132").Value = 1962.3334
$ws.Range("K132").Value = 5887.0002
$ws.Range("M132").Value = -3357.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14274.083
$ws.Range("I7").Value = 5682.8335
$ws.Range("J7").Value = 22865.334
$ws.Range("K7").Value = 5682.8335
$ws.Range("L7").Value = 22865.334
$ws.Range("M7").Value = -5570.8335
$ws.Range("N7").Value = -23089.334

$ws.Range("H40").Value = 3476668
$ws.Range("I40").Value = 4277591.5
$ws.Range("J40").Value = 5999.3335
$ws.Range("K40").Value = 4277591.5
$ws.Range("L40").Value = 5999.3335
$ws.Range("M40").Value = -4277455.5
$ws.Range("N40").Value = -6271.3335

$ws.Range("H94").Value = 59523.43

$ws.Range("H100").Value = 1733.1818
$ws.Range("I100").Value = 1631.6666
$ws.Range("J100").Value = 2190
$ws.Range("K100").Value = 1631.6666
$ws.Range("L100").Value = 2190
$ws.Range("M100").Value = -1090.6666
$ws.Range("N100").Value = -3272

$ws.Range("H122").Value = 3983.0476
$ws.Range("I122").Value = 2588
$ws.Range("K122").Value = 7764
$ws.Range("M122").Value = -5314

$ws.Range("H126").Value = 14274.083
$ws.Range("I126").Value = 5682.8335
$ws.Range("J126").Value = 22865.334
$ws.Range("K126").Value = 17048.5005
$ws.Range("L126").Value = 68596.00199999999
$ws.Range("M126").Value = -14578.5005
$ws.Range("N126").Value = -73536.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 46400.25
$ws.Range("J101").Value = 46400.25
$ws.Range("L101").Value = 46400.25
$ws.Range("N101").Value = -52890.25

$ws.Range("H122").Value = 1863.36
$ws.Range("I122").Value = 1878.5
$ws.Range("K122").Value = 5635.5
$ws.Range("M122").Value = -3185.5

$ws.Range("H126").Value = 2878.2593
$ws.Range("I126").Value = 2692.45
$ws.Range("J126").Value = 3409.1428
$ws.Range("K126").Value = 8077.349999999999
$ws.Range("L126").Value = 10227.4284
$ws.Range("M126").Value = -5607.349999999999
$ws.Range("N126").Value = -15167.4284

$ws.Range("H132").Value = 32875.562
$ws.Range("I132").Value = 46785.816
$ws.Range("J132").Value = 2273
$ws.Range("K132").Value = 140357.448
$ws.Range("L132").Value = 6819
$ws.Range("M132").Value = -137827.448
$ws.Range("N132").Value = -11879

$ws.Range("H136").Value = 21871.36
$ws.Range("I136").Value = 32314.938
$ws.Range("J136").Value = 3305
$ws.Range("K136").Value = 96944.814
$ws.Range("L136").Value = 9915
$ws.Range("M136").Value = -94394.814
$ws.Range("N136").Value = -15015
